$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from row 3 down into the new rows 4-6 so the
# newly added table rows share the same style as existing data rows.
$ws.Range("A3:K3").Copy()
$ws.Range("A4:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the table cells (rows 2-6) with the updated observation data.
$ws.Range("A2").Value = "us-core-average-blood-pressure"
$ws.Range("B2").Value = "US Core Average Blood Pressure Profile"
$ws.Range("C2").Value = "Observation Category Codes#vital-signs"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "LOINC#96607-7"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "dateTimeĵ, Periodĵ"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "US Core Average Blood Pressure Profile"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "LOINC#96608-5"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I3").Value = "optional"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "US Core Average Blood Pressure Profile"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "LOINC#96609-3"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I4").Value = "optional"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("A5").Value = "us-core-treatment-intervention-preference-bindings"
$ws.Range("B5").Value = "US Core Treatment Intervention Preference Bindings Profile"
$ws.Range("C5").Value = "US Core Category#treatment-intervention-preference"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "LOINC#75773-2"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H5").Value = "stringĵ, CodeableConceptĵ, Quantityĵ, booleanĵ, integerĵ, Rangeĵ, Ratioĵ, SampledDataĵ, timeĵ, dateTimeĵ, Periodĵ"
$ws.Range("I5").Value = "optional"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("A6").Value = "us-core-treatment-intervention-preference-grouping"
$ws.Range("B6").Value = "US Core Treatment Intervention Preference Grouping Profile"
$ws.Range("C6").Value = "US Core Category#treatment-intervention-preference"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "LOINC#75773-2"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H6").Value = "stringĵ, CodeableConceptĵ, Quantityĵ, booleanĵ, integerĵ, Rangeĵ, Ratioĵ, SampledDataĵ, timeĵ, dateTimeĵ, Periodĵ"
$ws.Range("I6").Value = "optional"
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
